# Add a new data row (row 14) to the time-registration sheet:
# Anna | 3 | 27.01.2021 | REA
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Anna"
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = "27.01.2021"
$ws.Range("D14").Value = "REA"

# Column D carries a column-level style (style id 3) that would otherwise
# get auto-applied to the new cell; the existing data rows (9-13) don't
# carry an explicit cell style, so reset D14 back to the workbook default
# to stay consistent with them.
$ws.Range("D14").Style = "Normal"
